$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsHPtFM = $wb.Worksheets.Item("HPtFM")

# "About" sheet: the old notes that explained "thermochemical water splitting"
# (rows 14-16) are no longer needed since that pathway row is being replaced,
# so remove those trailing rows.
$wsAbout.Rows("14:16").Delete()

# "HPtFM" sheet: rename the "thermochemical water splitting" pathway to
# "hydrocarbon partial oxidation" and mark it as using "heavy or residual
# fuel oil" (column I) instead of no fuel at all.
$wsHPtFM.Range("A6").Value = "hydrocarbon partial oxidation"
$wsHPtFM.Range("I6").Value = 1

# Restore cursor/selection state and active sheet to match the saved view:
# the About sheet cursor sits at F17, and the HPtFM sheet (now active/
# selected) has its cursor at I6.
$wsAbout.Range("F17").Select()
$wsHPtFM.Activate()
$wsHPtFM.Range("I6").Select()
